# Generate Report for Handoff
# Renames the two tracked localization items (new GUID-named files) and
# updates their status from "Handed back: in sync with en-US" to
# "Ready for handoff" across the Overview / zh-cn / de-de sheets, refreshes
# the per-language handoff file/datetime columns for the new handoff pass,
# and drops the now-unused "Latest Target File" / "Latest Handback File"
# columns (E/F) for the two rows since a handback hasn't happened yet.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "a9d6f234-fa1f-44ed-9bf7-c32a8516b059"
$newGuid1 = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d"
$oldGuid2 = "c4d713d3-6522-4898-b91a-9e7eaff7c0d7"
$newGuid2 = "ffff29a62014-0460-405d-88a9-e8168a4ad209"
$newHash  = "57064ced9f1784210164ffd6a2d120012a081521"

$statusNew = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: just the file names + status text change; layout is the same.
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid1.md"
$ov.Range("B2").Value = $statusNew
$ov.Range("C2").Value = $statusNew

$ov.Range("A3").Value = "$newGuid2.md"
$ov.Range("B3").Value = $statusNew
$ov.Range("C3").Value = $statusNew

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ecd36eb0b8b45043975de43632af85ed891ddc47/e2e/$newGuid1.md", [Type]::Missing, [Type]::Missing, "$newGuid1.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ecd36eb0b8b45043975de43632af85ed891ddc47/e2e/$newGuid2.md", [Type]::Missing, [Type]::Missing, "$newGuid2.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ecd36eb0b8b45043975de43632af85ed891ddc47/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de): same edits on both, only the
# handoff xlf file name / handoff datetime / hyperlink targets differ.
# ---------------------------------------------------------------------------
function Update-LangSheet($sheetName, $langCode, $handoffDatetime, $handoffUrlBase) {

    $ws = $wb.Worksheets.Item($sheetName)
    $xlf = "$newGuid1.$newHash.$langCode.xlf"

    # Row 2 (first tracked item)
    $ws.Range("A2").Value = "$newGuid1.md"
    $ws.Range("B2").Value = $statusNew
    $ws.Range("C2").Value = $xlf
    $ws.Range("D2").Value = $handoffDatetime
    $ws.Range("E2").ClearContents()
    $ws.Range("F2").ClearContents()
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Include"

    # Row 3 (second tracked item - now shares the same handoff file as row 2)
    $ws.Range("A3").Value = "$newGuid2.md"
    $ws.Range("B3").Value = $statusNew
    $ws.Range("C3").Value = $xlf
    $ws.Range("D3").Value = $handoffDatetime
    $ws.Range("E3").ClearContents()
    $ws.Range("F3").ClearContents()
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Include"

    # Row 4 (.localization-config, not localized) - only the two date-ish
    # cells actually change (they reuse the "not handed back yet" sentinel).
    $ws.Range("D4").Value = "0001-01-01 00:00:00"
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Ignored"

    # Rebuild the hyperlinks cleanly: drop the stale Latest Target
    # File / Latest Handback File links (E2/F2/E3/F3) entirely, and point
    # the remaining ones at the new file names.
    $ws.Hyperlinks.Delete()

    $mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/ecd36eb0b8b45043975de43632af85ed891ddc47/e2e/$newGuid1.md"
    $mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/ecd36eb0b8b45043975de43632af85ed891ddc47/e2e/$newGuid2.md"
    $xlfUrl = "$handoffUrlBase/$xlf"
    $cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ecd36eb0b8b45043975de43632af85ed891ddc47/.localization-config"

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, "$newGuid1.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl, [Type]::Missing, [Type]::Missing, $xlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, "$newGuid2.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $xlfUrl, [Type]::Missing, [Type]::Missing, $xlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
}

$zhHandoffUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e08d96c7979e78334b84335f563951b7b90af62/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/425c93f6d57dc3d9b15bd481412bddac08dd58c9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

Update-LangSheet "zh-cn" "zh-cn" "2016-03-11 01:32:56" $zhHandoffUrlBase
Update-LangSheet "de-de" "de-de" "2016-03-11 01:33:03" $deHandoffUrlBase
